$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Update the "Ende" time for C7 from 16:00 to 16:30 (stored as fraction of a day)
$ws.Range("C7").Value = 0.6875

# Update selection / active cell shown in the sheet view from J10 to J8
$ws.Range("J8").Select()

$wb.Save()
